$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 117
$ws1.Range("F3").Value = 535
$ws1.Range("F4").Value = 1538
$ws1.Range("F5").Value = 155
$ws1.Range("F7").Value = 5211
$ws1.Range("F8").Value = 166
$ws1.Range("F9").Value = 750
$ws1.Range("F12").Value = 333
$ws1.Range("F13").Value = 57
$ws1.Range("F14").Value = 9
$ws1.Range("F15").Value = 6459
$ws1.Range("F16").Value = 19
$ws1.Range("F19").Value = 158
$ws1.Range("F21").Value = 15435
$ws1.Range("F22").Value = 1528
$ws1.Range("F23").Value = 288
$ws1.Range("F24").Value = 143
$ws1.Range("F25").Value = 102
$ws1.Range("F26").Value = 11071
$ws1.Range("F27").Value = 758
$ws1.Range("F28").Value = 4334
$ws1.Range("F29").Value = 245
$ws1.Range("F32").Value = 305

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 117
$ws4.Range("F3").Value = 535
$ws4.Range("F4").Value = 1538
$ws4.Range("F5").Value = 155
$ws4.Range("F8").Value = 5211
$ws4.Range("F9").Value = 166
$ws4.Range("F10").Value = 750
$ws4.Range("F14").Value = 333
$ws4.Range("F15").Value = 57
$ws4.Range("F17").Value = 9
$ws4.Range("F18").Value = 6459
$ws4.Range("F19").Value = 19
$ws4.Range("F22").Value = 158
$ws4.Range("F24").Value = 15435
$ws4.Range("F25").Value = 1528
$ws4.Range("F26").Value = 288
$ws4.Range("F27").Value = 143
$ws4.Range("F28").Value = 102
$ws4.Range("F29").Value = 11071
$ws4.Range("F30").Value = 758
$ws4.Range("F31").Value = 4334
$ws4.Range("F32").Value = 245
$ws4.Range("F35").Value = 305
